$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Figura" in column O (next to TipoContrato), reflecting
# recognition of child cases ("Reconocimiento de casos hijos agregado")
$ws.Range("O1").Value = "Figura"

# Update the sheet's current selection to the new working range
$excel.Goto($ws.Range("A2:W7"))
